$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-17 05:48:31"
$ws.Range("M2").Value = "1.3 °C 5:05 TU"
$ws.Range("E3").Value = "2026-02-17 05:48:33"
$ws.Range("O3").Value = "-7.1 °C"
$ws.Range("E4").Value = "2026-02-17 05:48:36"
$ws.Range("J4").Value = "1015.4 hPa"
$ws.Range("E5").Value = "2026-02-17 05:48:38"
$ws.Range("O5").Value = "-6.5 °C"
$ws.Range("E6").Value = "2026-02-17 05:48:40"
$ws.Range("H6").Value = "'88%"
$ws.Range("J6").Value = "1015.1 hPa"
$ws.Range("E7").Value = "2026-02-17 05:48:43"
$ws.Range("J7").Value = "1014.7 hPa"
$ws.Range("N7").Value = "12.6 °C 5:22 TU"
$ws.Range("O7").Value = "14.1 °C"
$ws.Range("E8").Value = "2026-02-17 05:48:45"
$ws.Range("J8").Value = "1015.0 hPa"
$ws.Range("N8").Value = "8.7 °C 5:11 TU"
$ws.Range("O8").Value = "9.9 °C"
$ws.Range("E9").Value = "2026-02-17 05:48:47"
$ws.Range("O9").Value = "12.4 °C"
$ws.Range("E10").Value = "2026-02-17 05:48:50"
$ws.Range("H10").Value = "'88%"
$ws.Range("E11").Value = "2026-02-17 05:48:52"
$ws.Range("H11").Value = "'38%"
$ws.Range("N11").Value = "1.9 °C 5:24 TU"
$ws.Range("O11").Value = "6.3 °C"
$ws.Range("E12").Value = "2026-02-17 05:48:54"
$ws.Range("H12").Value = "'47%"
$ws.Range("N12").Value = "10.9 °C 5:27 TU"
$ws.Range("O12").Value = "12.8 °C"
$ws.Range("E13").Value = "2026-02-17 05:48:56"
$ws.Range("J13").Value = "1016.8 hPa"
$ws.Range("O13").Value = "4.0 °C"
$ws.Range("E14").Value = "2026-02-17 05:48:59"
$ws.Range("N14").Value = "10.7 °C 5:29 TU"
$ws.Range("O14").Value = "13.2 °C"
$ws.Range("E15").Value = "2026-02-17 05:49:01"
$ws.Range("H15").Value = "'47%"
$ws.Range("N15").Value = "10.6 °C 5:29 TU"
$ws.Range("O15").Value = "12.3 °C"
$ws.Range("E16").Value = "2026-02-17 05:49:03"
$ws.Range("H16").Value = "'47%"
$ws.Range("E17").Value = "2026-02-17 05:49:06"
$ws.Range("H17").Value = "'50%"
$ws.Range("E18").Value = "2026-02-17 05:49:08"
$ws.Range("J18").Value = "1015.4 hPa"
$ws.Range("M18").Value = "6.9 °C 5:29 TU"
$ws.Range("E19").Value = "2026-02-17 05:49:10"
$ws.Range("H19").Value = "'68%"
$ws.Range("K19").Value = "-0.1 MJ/m2"
$ws.Range("E20").Value = "2026-02-17 05:49:12"
$ws.Range("H20").Value = "'38%"
$ws.Range("E21").Value = "2026-02-17 05:49:15"
$ws.Range("H21").Value = "'29%"
$ws.Range("J21").Value = "1015.7 hPa"
$ws.Range("O21").Value = "7.7 °C"
$ws.Range("E22").Value = "2026-02-17 05:49:17"
$ws.Range("E23").Value = "2026-02-17 05:49:19"
$ws.Range("E24").Value = "2026-02-17 05:49:22"
$ws.Range("H24").Value = "'74%"
$ws.Range("J24").Value = "1017.5 hPa"
$ws.Range("N24").Value = "9.5 °C 5:23 TU"
$ws.Range("E25").Value = "2026-02-17 05:49:24"
$ws.Range("E26").Value = "2026-02-17 05:49:26"
$ws.Range("E27").Value = "2026-02-17 05:49:29"
$ws.Range("H27").Value = "'41%"
$ws.Range("N27").Value = "-3.8 °C 5:19 TU"
$ws.Range("O27").Value = "-2.7 °C"
$ws.Range("E28").Value = "2026-02-17 05:49:31"
$ws.Range("J28").Value = "1015.6 hPa"
$ws.Range("E29").Value = "2026-02-17 05:49:33"
$ws.Range("O29").Value = "11.7 °C"
$ws.Range("E30").Value = "2026-02-17 05:49:36"
$ws.Range("H30").Value = "'47%"
$ws.Range("J30").Value = "1014.6 hPa"
$ws.Range("O30").Value = "11.9 °C"
$ws.Range("E31").Value = "2026-02-17 05:49:38"
$ws.Range("H31").Value = "'66%"
$ws.Range("J31").Value = "1015.3 hPa"
$ws.Range("N31").Value = "7.7 °C 5:07 TU"
$ws.Range("O31").Value = "9.7 °C"
$ws.Range("E32").Value = "2026-02-17 05:49:40"
$ws.Range("E33").Value = "2026-02-17 05:49:42"
$ws.Range("H33").Value = "'38%"
$ws.Range("K33").Value = "-0.1 MJ/m2"
$ws.Range("L33").Value = "51.5 km/h - 291º 5:03 TU"
$ws.Range("O33").Value = "4.7 °C"
$ws.Range("E34").Value = "2026-02-17 05:49:45"
$ws.Range("E35").Value = "2026-02-17 05:49:47"
$ws.Range("I35").Value = "2.5 mm"
$ws.Range("J35").Value = "1018.2 hPa"
$ws.Range("O35").Value = "5.4 °C"
$ws.Range("E36").Value = "2026-02-17 05:49:50"
$ws.Range("H36").Value = "'50%"
$ws.Range("J36").Value = "1015.4 hPa"
$ws.Range("N36").Value = "10.2 °C 5:20 TU"
$ws.Range("O36").Value = "12.7 °C"
$ws.Range("E37").Value = "2026-02-17 05:49:52"
$ws.Range("H37").Value = "'46%"
$ws.Range("J37").Value = "1015.6 hPa"
$ws.Range("N37").Value = "4.6 °C 5:21 TU"
$ws.Range("O37").Value = "7.9 °C"
$ws.Range("E38").Value = "2026-02-17 05:49:54"
$ws.Range("E39").Value = "2026-02-17 05:49:57"
$ws.Range("H39").Value = "'53%"
$ws.Range("I39").Value = "1.1 mm"
$ws.Range("M39").Value = "-2.8 °C 5:20 TU"
$ws.Range("O39").Value = "-4.7 °C"
$ws.Range("E40").Value = "2026-02-17 05:49:59"
$ws.Range("J40").Value = "1017.4 hPa"
$ws.Range("E41").Value = "2026-02-17 05:50:01"
$ws.Range("H41").Value = "'49%"
$ws.Range("J41").Value = "1015.4 hPa"
$ws.Range("N41").Value = "12.0 °C 5:05 TU"
$ws.Range("O41").Value = "14.9 °C"
$ws.Range("E42").Value = "2026-02-17 05:50:04"
$ws.Range("H42").Value = "'48%"
$ws.Range("O42").Value = "12.8 °C"
$ws.Range("E43").Value = "2026-02-17 05:50:06"
$ws.Range("N43").Value = "2.1 °C 5:29 TU"
$ws.Range("O43").Value = "4.3 °C"
$ws.Range("E44").Value = "2026-02-17 05:50:08"
$ws.Range("M44").Value = "-4.9 °C 5:27 TU"
$ws.Range("O44").Value = "-5.4 °C"
$ws.Range("E45").Value = "2026-02-17 05:50:11"
$ws.Range("L45").Value = "39.2 km/h - 206º 5:10 TU"
$ws.Range("E46").Value = "2026-02-17 05:50:13"
$ws.Range("K46").Value = "-0.1 MJ/m2"
$ws.Range("O46").Value = "13.4 °C"
